$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shift the "7.0 / Pagina de registro" and "8.0 / Conocer mas sobre ti"
# rows down by one (row 20->21, row 21->22) to make room for the new
# "6.1 / Modificacion de contrasena" row, working bottom-up so we don't
# clobber data we still need to read.
$ws.Range("A22").Value = $ws.Range("A21").Value()
$ws.Range("B22").Value = $ws.Range("B21").Value()

$ws.Range("A21").Value = $ws.Range("A20").Value()
$ws.Range("B21").Value = $ws.Range("B20").Value()

# --- Fill the newly freed row 20 with the new entry. Column A holds a
# version-number-looking string ("6.1") that must stay text (matching the
# rest of the column), not be auto-converted to a number, so force the
# text number format before assigning it and then drop back to the
# sheet's normal style so no stray formatting is left behind.
$ws.Range("B20").Value = "Modificación de contraseña"
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "6.1"
$ws.Range("A20").Style = "Normal"

# --- Restore the view state recorded in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B18").Select()
